$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.009.01"
$ws.Range("E2").Value = "  +1.37%  "

$ws.Range("D3").Value = "1.978.61"
$ws.Range("E3").Value = "  +1.05%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.21%  "

$ws.Range("E6").Value = "  +1.86%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.89"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.80%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  +0.98%  "

$ws.Range("E10").Value = "  -1.64%  "

$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.839"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.16%  "

$ws.Range("D15").Value = "2.270.05"
$ws.Range("E15").Value = "  +1.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.70%  "

$ws.Range("D17").Value = "1.981.28"
$ws.Range("E17").Value = "  +1.14%  "

$ws.Range("D18").Value = "36.864.05"
$ws.Range("E18").Value = "  +1.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.44%  "

$ws.Range("D20").Value = "0.0₃0857"
$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.73%  "

$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.38%  "

$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.146"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.60%  "

$ws.Range("E27").Value = "  +0.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.67%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.79%  "

$ws.Range("E30").Value = "  +16.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.121"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0620"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("E34").Value = "  +5.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.83%  "

$ws.Range("E36").Value = "  +0.19%  "

$ws.Range("E37").Value = "  +1.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.97%  "

$ws.Range("E39").Value = "  -5.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0976"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.28%  "

$ws.Range("E41").Value = "  +1.34%  "

$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0212"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.47%  "

$ws.Range("E44").Value = "  +3.29%  "

$ws.Range("D45").Value = "1.366.04"
$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.23%  "

$ws.Range("E47").Value = "  -0.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.05%  "

$ws.Range("E51").Value = "  +10.17%  "
